$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the leaderboard data (rows 10-12, columns A-D) leaving styles intact
$ws.Range("A10:D12").ClearContents()

# Update the selection to match the new state (A12:A13, active cell on A13)
$ws.Range("A12:A13").Select()
